$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 209.875
$ws.Cells.Item(4, 9).Value = 209.875
$ws.Cells.Item(4, 11).Value = 209.875
$ws.Cells.Item(4, 13).Value = -95.875

$ws.Cells.Item(8, 8).Value = 719.8333
$ws.Cells.Item(8, 9).Value = 70.888885
$ws.Cells.Item(8, 10).Value = 2666.6667
$ws.Cells.Item(8, 11).Value = 212.666655
$ws.Cells.Item(8, 12).Value = 8000.000100000001
$ws.Cells.Item(8, 13).Value = -73.66665499999999
$ws.Cells.Item(8, 14).Value = -8278.000100000001

$ws.Cells.Item(18, 8).Value = 1547863.4
$ws.Cells.Item(18, 9).Value = 1985824
$ws.Cells.Item(18, 10).Value = 15001
$ws.Cells.Item(18, 11).Value = 1985824
$ws.Cells.Item(18, 12).Value = 15001
$ws.Cells.Item(18, 13).Value = -1985540
$ws.Cells.Item(18, 14).Value = -15569

$ws.Cells.Item(53, 8).Value = 320
$ws.Cells.Item(53, 9).Value = 255.84616
$ws.Cells.Item(53, 10).Value = 372.125
$ws.Cells.Item(53, 11).Value = 255.84616
$ws.Cells.Item(53, 12).Value = 372.125
$ws.Cells.Item(53, 13).Value = 381.15384
$ws.Cells.Item(53, 14).Value = -1646.125

$ws.Cells.Item(137, 8).Value = 3849673.8
$ws.Cells.Item(137, 9).Value = 10005950
$ws.Cells.Item(137, 10).Value = 2001.0625
$ws.Cells.Item(137, 11).Value = 30017850
$ws.Cells.Item(137, 12).Value = 6003.1875
$ws.Cells.Item(137, 13).Value = -30015300
$ws.Cells.Item(137, 14).Value = -11103.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2191.02
$ws.Cells.Item(32, 9).Value = 1397.8068
$ws.Cells.Item(32, 10).Value = 8007.9165
$ws.Cells.Item(32, 11).Value = 1397.8068
$ws.Cells.Item(32, 12).Value = 8007.9165
$ws.Cells.Item(32, 13).Value = -1110.8068
$ws.Cells.Item(32, 14).Value = -8581.916499999999

$ws.Cells.Item(61, 8).Value = 2752.7368
$ws.Cells.Item(61, 9).Value = 1123.8096
$ws.Cells.Item(61, 10).Value = 4764.9414
$ws.Cells.Item(61, 11).Value = 1123.8096
$ws.Cells.Item(61, 12).Value = 4764.9414
$ws.Cells.Item(61, 13).Value = -911.8096
$ws.Cells.Item(61, 14).Value = -5188.9414

$ws.Cells.Item(122, 8).Value = 2927.5789
$ws.Cells.Item(122, 9).Value = 1829.2858
$ws.Cells.Item(122, 10).Value = 6002.8
$ws.Cells.Item(122, 11).Value = 5487.857400000001
$ws.Cells.Item(122, 12).Value = 18008.4
$ws.Cells.Item(122, 13).Value = -3037.857400000001
$ws.Cells.Item(122, 14).Value = -22908.4

$ws.Cells.Item(132, 8).Value = 14288847
$ws.Cells.Item(132, 9).Value = 20003172
$ws.Cells.Item(132, 10).Value = 3034.9
$ws.Cells.Item(132, 11).Value = 60009516
$ws.Cells.Item(132, 12).Value = 9104.700000000001
$ws.Cells.Item(132, 13).Value = -60006986
$ws.Cells.Item(132, 14).Value = -14164.7

$ws.Cells.Item(136, 8).Value = 2752.7368
$ws.Cells.Item(136, 9).Value = 1123.8096
$ws.Cells.Item(136, 10).Value = 4764.9414
$ws.Cells.Item(136, 11).Value = 3371.4288
$ws.Cells.Item(136, 12).Value = 14294.8242
$ws.Cells.Item(136, 13).Value = -821.4288000000001
$ws.Cells.Item(136, 14).Value = -19394.8242

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1657.15
$ws.Cells.Item(105, 9).Value = 1356.1538
$ws.Cells.Item(105, 10).Value = 2216.1428
$ws.Cells.Item(105, 11).Value = 1356.1538
$ws.Cells.Item(105, 12).Value = 2216.1428
$ws.Cells.Item(105, 13).Value = 390.8462
$ws.Cells.Item(105, 14).Value = -5710.1428

$ws.Cells.Item(134, 8).Value = 3579.1353
$ws.Cells.Item(134, 9).Value = 3379.5715
$ws.Cells.Item(134, 10).Value = 4200
$ws.Cells.Item(134, 11).Value = 10138.7145
$ws.Cells.Item(134, 12).Value = 12600
$ws.Cells.Item(134, 13).Value = -7603.7145
$ws.Cells.Item(134, 14).Value = -17670

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1451451
$ws.Cells.Item(31, 9).Value = 2501354.5
$ws.Cells.Item(31, 10).Value = 3308.3794
$ws.Cells.Item(31, 11).Value = 2501354.5
$ws.Cells.Item(31, 12).Value = 3308.3794
$ws.Cells.Item(31, 13).Value = -2501059.5
$ws.Cells.Item(31, 14).Value = -3898.3794

$ws.Cells.Item(34, 8).Value = 1451451
$ws.Cells.Item(34, 9).Value = 2501354.5
$ws.Cells.Item(34, 10).Value = 3308.3794
$ws.Cells.Item(34, 11).Value = 2501354.5
$ws.Cells.Item(34, 12).Value = 3308.3794
$ws.Cells.Item(34, 13).Value = -2501152.5
$ws.Cells.Item(34, 14).Value = -3712.3794

$ws.Cells.Item(99, 8).Value = 2073.818
$ws.Cells.Item(99, 9).Value = 1171.2858
$ws.Cells.Item(99, 10).Value = 3653.25
$ws.Cells.Item(99, 11).Value = 1171.2858
$ws.Cells.Item(99, 12).Value = 3653.25
$ws.Cells.Item(99, 13).Value = 326.7141999999999
$ws.Cells.Item(99, 14).Value = -6649.25

$ws.Cells.Item(122, 8).Value = 3892.9285
$ws.Cells.Item(122, 9).Value = 3543.1428
$ws.Cells.Item(122, 10).Value = 4242.7144
$ws.Cells.Item(122, 11).Value = 10629.4284
$ws.Cells.Item(122, 12).Value = 12728.1432
$ws.Cells.Item(122, 13).Value = -8179.428400000001
$ws.Cells.Item(122, 14).Value = -17628.1432

$ws.Cells.Item(126, 8).Value = 2073.818
$ws.Cells.Item(126, 9).Value = 1171.2858
$ws.Cells.Item(126, 10).Value = 3653.25
$ws.Cells.Item(126, 11).Value = 3513.8574
$ws.Cells.Item(126, 12).Value = 10959.75
$ws.Cells.Item(126, 13).Value = -1043.8574
$ws.Cells.Item(126, 14).Value = -15899.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(130, 8).Value = 3000
$ws.Cells.Item(130, 9).Value = 5000
$ws.Cells.Item(130, 10).Value = 2500
$ws.Cells.Item(130, 11).Value = 15000
$ws.Cells.Item(130, 12).Value = 7500
$ws.Cells.Item(130, 13).Value = -9980
$ws.Cells.Item(130, 14).Value = -17540

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1079.5667
$ws.Cells.Item(97, 9).Value = 801.5217
$ws.Cells.Item(97, 11).Value = 801.5217
$ws.Cells.Item(97, 13).Value = -305.5217

$ws.Cells.Item(102, 8).Value = 2089.5217
$ws.Cells.Item(102, 9).Value = 1289.9286
$ws.Cells.Item(102, 10).Value = 3333.3333
$ws.Cells.Item(102, 11).Value = 1289.9286
$ws.Cells.Item(102, 12).Value = 3333.3333
$ws.Cells.Item(102, 13).Value = 332.0714
$ws.Cells.Item(102, 14).Value = -6577.3333

$ws.Cells.Item(132, 8).Value = 3827.9062
$ws.Cells.Item(132, 9).Value = 2699.9
$ws.Cells.Item(132, 10).Value = 5707.9165
$ws.Cells.Item(132, 11).Value = 8099.700000000001
$ws.Cells.Item(132, 12).Value = 17123.7495
$ws.Cells.Item(132, 13).Value = -5569.700000000001
$ws.Cells.Item(132, 14).Value = -22183.7495

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2630.8
$ws.Cells.Item(7, 9).Value = 1704
$ws.Cells.Item(7, 10).Value = 2862.5
$ws.Cells.Item(7, 11).Value = 1704
$ws.Cells.Item(7, 12).Value = 2862.5
$ws.Cells.Item(7, 13).Value = -1592
$ws.Cells.Item(7, 14).Value = -3086.5

$ws.Cells.Item(40, 8).Value = 3317.1667
$ws.Cells.Item(40, 9).Value = 1999.5
$ws.Cells.Item(40, 10).Value = 3976
$ws.Cells.Item(40, 11).Value = 1999.5
$ws.Cells.Item(40, 12).Value = 3976
$ws.Cells.Item(40, 13).Value = -1863.5
$ws.Cells.Item(40, 14).Value = -4248

$ws.Cells.Item(126, 8).Value = 2630.8
$ws.Cells.Item(126, 9).Value = 1704
$ws.Cells.Item(126, 10).Value = 2862.5
$ws.Cells.Item(126, 11).Value = 5112
$ws.Cells.Item(126, 12).Value = 8587.5
$ws.Cells.Item(126, 13).Value = -2642
$ws.Cells.Item(126, 14).Value = -13527.5

$ws.Cells.Item(127, 8).Value = 28000
$ws.Cells.Item(127, 10).Value = 28000
$ws.Cells.Item(127, 12).Value = 28000
$ws.Cells.Item(127, 14).Value = -37920

$ws.Cells.Item(132, 8).Value = 2334.9714
$ws.Cells.Item(132, 9).Value = 1305.5
$ws.Cells.Item(132, 11).Value = 3916.5
$ws.Cells.Item(132, 13).Value = -1386.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 3197313
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 3197313
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 3197313
$ws.Cells.Item(2, 14).Value = -3197537
$ws.Cells.Item(2, 13).ClearContents()

$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()

$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 2082.8
$ws.Cells.Item(113, 9).Value = 570
$ws.Cells.Item(113, 10).Value = 3595.6
$ws.Cells.Item(113, 11).Value = 1710
$ws.Cells.Item(113, 12).Value = 10786.8
$ws.Cells.Item(113, 13).Value = 460
$ws.Cells.Item(113, 14).Value = -15126.8

$ws.Cells.Item(132, 8).Value = 243112.23
$ws.Cells.Item(132, 9).Value = 358972
$ws.Cells.Item(132, 10).Value = 11392.714
$ws.Cells.Item(132, 11).Value = 1076916
$ws.Cells.Item(132, 12).Value = 34178.142
$ws.Cells.Item(132, 13).Value = -1074386
$ws.Cells.Item(132, 14).Value = -39238.142
